{"js": "// The \"1. fatty acids (noun)\" entry is wrong/misspelled and its two\n// definition lines need to go away entirely; it becomes a single\n// \"1. awef     (missing)\" paragraph (matching the other \"missing\" entries\n// elsewhere in the glossary), and the two following definition paragraphs\n// are removed outright.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst headingIndex = items.findIndex((p) => p.text.trim() === \"1. fatty acids     (noun)\");\n\nif (headingIndex === -1) {\n  throw new Error('Could not find paragraph \"1. fatty acids     (noun)\"');\n}\n\n// Replace the heading text in place.\nitems[headingIndex].insertText(\"1. awef     (missing)\", Word.InsertLocation.replace);\n\n// Remove the two definition paragraphs that followed the heading.\nconst defs = [];\nfor (let i = headingIndex + 1; i < items.length && i < headingIndex + 3; i++) {\n  const t = items[i].text.trim();\n  if (t.startsWith(\"-\")) {\n    defs.push(items[i]);\n  }\n}\ndefs.forEach((p) => p.delete());\n\nawait context.sync();\n", "ps1": "# The \"1. fatty acids (noun)\" glossary entry is wrong/misspelled and its\n# two definition lines need to be removed entirely; it becomes a single\n# \"1. awef     (missing)\" paragraph, and the two following definition\n# paragraphs (the \"-  any of ...\" lines) are deleted outright.\n\n$d = $word.ActiveDocument\n\n# Locate the heading paragraph by its current text.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13) -eq \"1. fatty acids     (noun)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find paragraph '1. fatty acids     (noun)'\"\n}\n\n# Grab the two definition paragraphs that immediately follow the heading\n# before we change the heading text.\n$def1 = $target.Next()\n$def2 = $def1.Next()\n\n# Fix the misspelled/missing heading text (keeps the paragraph mark).\n$target.Range.Text = \"1. awef     (missing)\"\n\n# Remove the two now-orphaned definition paragraphs.\n$delRange = $d.Range($def1.Range.Start, $def2.Range.End)\n$delRange.Delete()\n"}
